# USER PHP TESTING.xlsx - "Added all testing doc files"
#
# This script:
#   1. Fixes a pre-existing C36/C37 data mismatch and adds a new C38 value
#      (Check reviews / Check product details / Write reviews rows).
#   2. Inserts new test-case rows for a "Checkout" -> PayPal flow
#      ("Proceed to Checkout", "Collection Slot") between the existing
#      "Cart" and "Payment" sections.
#   3. Inserts new test-case rows under "Payment" for "PayPal Login",
#      "Paypal Payment" and "Invoice".
#   4. Updates the window view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up the "Product" section (reviews) rows -----------------------
# C36 previously pointed at the "All the product details..." text (meant
# for row 37); it should hold the (new) "user can check all the reviews"
# text. C37 should hold the "All the product details..." text that used to
# sit in C36. A new C38 cell documents what "Write reviews" asks the user.
$ws.Range("C36").Value = "User can check all the review  a particular product has"
$ws.Range("C37").Value = "All the product details such as allergen information , price and stock are displayed"
$ws.Range("C38").Value = "Asks the user to write their review and provide ratings in the form of stars"

# --- 2. Make room for the new rows -----------------------------------------
# Old layout (rows 44-57):
#   44 Checkout (header)               51 About Us (header)
#   45 (blank / wrap style only)       52 About Us content row
#   46 Payment (header)                55 Contact Us (header)
#                                      56-57 Contact Us content rows
#
# New layout (rows 45-60):
#   45 Checkout (header)               54 About Us (header)
#   46 Proceed to Checkout             55 About Us content row
#   47 Collection Slot                 58 Contact Us (header)
#   48 (blank / wrap style only)      59-60 Contact Us content rows
#   49 Payment (header)
#   50 PayPal Login
#   51 Paypal Payment
#   52 Invoice
#
# Net effect: 1 row inserted right before old row 44 (Checkout header), and
# 2 more rows inserted right before old row 46 (Payment header). Doing the
# lower insert first keeps the row44 address stable for the second insert.
$ws.Range("A46:A47").EntireRow.Insert()
$ws.Rows(44).Insert()

# --- 3. Checkout section (rows 45-48) --------------------------------------
$ws.Range("A46").Value = "Proceed to Checkout"
$ws.Range("B46").Value = "Click on the checkout button "
$ws.Range("C46").Value = "Once the user checksout, they are directed to Paypal login"
$ws.Range("D46").Value = "Pass"

$ws.Range("A47").Value = "Collection Slot"
$ws.Range("D47").Value = "Pass"

# The old "Checkout" header row carried a stray D44="Pass" and no C-style
# wrap cell beyond C44; after the insert it now sits at row 45 and must not
# carry a D value any more (the new layout keeps only A45 + a styled, empty
# C45).
$ws.Range("D45").ClearContents()

# --- 4. Payment section (rows 49-52) ---------------------------------------
$ws.Range("A50").Value = "PayPal Login"
$ws.Range("B50").Value = "User is directed to paypal login after proceeding to checkout"
$ws.Range("C50").Value = "User is asked to enter their paypal details"
$ws.Range("D50").Value = "Pass"

$ws.Range("A51").Value = "Paypal Payment"
$ws.Range("B51").Value = "User is directeed to payment after logging in to Paypal"
$ws.Range("C51").Value = "User is prompted to complete the purchase"
$ws.Range("D51").Value = "Pass"

$ws.Range("A52").Value = "Invoice"
$ws.Range("B52").Value = "Once the payment is done they will receive an invoice"
$ws.Range("C52").Value = "The user will receive their invoice in their mail address"
$ws.Range("D52").Value = "Pass"

# --- 5. View state -----------------------------------------------------------
$excel.ActiveWindow.Zoom = 55
$ws.Range("D45").Select()
